$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 840
$ws.Range("I5").Value = 675
$ws.Range("K5").Value = 675
$ws.Range("M5").Value = -560
# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 7239
$ws.Range("J32").Value = 8361.25
$ws.Range("L32").Value = 8361.25
$ws.Range("N32").Value = -9013.25
# Row 99 (Leve Item ID 19883)
$ws.Range("H99").Value = 9665
$ws.Range("I99").Value = 9497.5
$ws.Range("K99").Value = 28492.5
$ws.Range("M99").Value = -26994.5
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 34385200
$ws.Range("I116").Value = 20930256
$ws.Range("J116").Value = 66677068
$ws.Range("K116").Value = 20930256
$ws.Range("L116").Value = 66677068
$ws.Range("M116").Value = -20926814
$ws.Range("N116").Value = -66683952
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 37040716
$ws.Range("I137").Value = 62501850
$ws.Range("J137").Value = 6340.4546
$ws.Range("K137").Value = 187505550
$ws.Range("L137").Value = 19021.3638
$ws.Range("M137").Value = -187503000
$ws.Range("N137").Value = -24121.3638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 16669281
$ws.Range("J88").Value = 2904.111
$ws.Range("L88").Value = 2904.111
$ws.Range("N88").Value = -3716.111
# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 16669281
$ws.Range("J91").Value = 2904.111
$ws.Range("L91").Value = 2904.111
$ws.Range("N91").Value = -5712.111

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 29413966
$ws.Range("I86").Value = 31252026
$ws.Range("J86").Value = 5006
$ws.Range("K86").Value = 31252026
$ws.Range("L86").Value = 5006
$ws.Range("M86").Value = -31250903
$ws.Range("N86").Value = -7252
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 29413966
$ws.Range("I89").Value = 31252026
$ws.Range("J89").Value = 5006
$ws.Range("K89").Value = 156260130
$ws.Range("L89").Value = 25030
$ws.Range("M89").Value = -156254514
$ws.Range("N89").Value = -36262
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 27800894
$ws.Range("I107").Value = 17554.182
$ws.Range("J107").Value = 71460424
$ws.Range("K107").Value = 17554.182
$ws.Range("L107").Value = 71460424
$ws.Range("M107").Value = -15634.182
$ws.Range("N107").Value = -71464264
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2616.0667
$ws.Range("I134").Value = 2267.2144
$ws.Range("K134").Value = 6801.6432
$ws.Range("M134").Value = -4266.6432

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 2322.3572
$ws.Range("I16").Value = 1810.875
$ws.Range("K16").Value = 1810.875
$ws.Range("M16").Value = -1523.875
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 648
$ws.Range("I107").Value = 658.2778
$ws.Range("K107").Value = 658.2778
$ws.Range("M107").Value = 1261.7222
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 2322.3572
$ws.Range("I113").Value = 1810.875
$ws.Range("K113").Value = 1810.875
$ws.Range("M113").Value = 359.125
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3327.7083
$ws.Range("J134").Value = 3374.5
$ws.Range("L134").Value = 10123.5
$ws.Range("N134").Value = -15193.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 61 (Leve Item ID 4727)
$ws.Range("H61").Value = 343.4
$ws.Range("I61").Value = 79.25
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 237.75
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -22.75
$ws.Range("N61").Value = -4630
# Row 101 (Leve Item ID 19820)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 100769.5
$ws.Range("I121").Value = 20314
$ws.Range("J121").Value = 158237.72
$ws.Range("K121").Value = 60942
$ws.Range("L121").Value = 474713.16
$ws.Range("M121").Value = -59632
$ws.Range("N121").Value = -477333.16
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 7501.3613
$ws.Range("I131").Value = 1031.0834
$ws.Range("J131").Value = 10736.5
$ws.Range("K131").Value = 3093.2502
$ws.Range("L131").Value = 32209.5
$ws.Range("M131").Value = 1946.7498
$ws.Range("N131").Value = -42289.5
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1340.1
$ws.Range("I132").Value = 1175.5
$ws.Range("J132").Value = 1998.5
$ws.Range("K132").Value = 10579.5
$ws.Range("L132").Value = 17986.5
$ws.Range("M132").Value = -8049.5
$ws.Range("N132").Value = -23046.5
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 1673.625
$ws.Range("I140").Value = 1364.8334
$ws.Range("J140").Value = 2600
$ws.Range("K140").Value = 4094.5002
$ws.Range("L140").Value = 7800
$ws.Range("M140").Value = 1085.4998
$ws.Range("N140").Value = -18160

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 1888.1875
$ws.Range("I97").Value = 1916.2307
$ws.Range("K97").Value = 1916.2307
$ws.Range("M97").Value = -1420.2307
# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 715.35
$ws.Range("I107").Value = 538
$ws.Range("J107").Value = 1424.75
$ws.Range("K107").Value = 538
$ws.Range("L107").Value = 1424.75
$ws.Range("M107").Value = 1382
$ws.Range("N107").Value = -5264.75
# Row 139 (Leve Item ID 42373)
$ws.Range("H139").Value = 64210.527
$ws.Range("J139").Value = 64210.527
$ws.Range("L139").Value = 64210.527
$ws.Range("N139").Value = -74490.527

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2482311.8
$ws.Range("I22").Value = 883.6667
$ws.Range("K22").Value = 883.6667
$ws.Range("M22").Value = -588.6667
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2482311.8
$ws.Range("I27").Value = 883.6667
$ws.Range("K27").Value = 883.6667
$ws.Range("M27").Value = -776.6667
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2319.6897
$ws.Range("I46").Value = 1817.6364
$ws.Range("J46").Value = 2626.5
$ws.Range("K46").Value = 1817.6364
$ws.Range("L46").Value = 2626.5
$ws.Range("M46").Value = -1629.6364
$ws.Range("N46").Value = -3002.5
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 432.16666
$ws.Range("I55").Value = 499.6154
$ws.Range("J55").Value = 352.45456
$ws.Range("K55").Value = 499.6154
$ws.Range("L55").Value = 352.45456
$ws.Range("M55").Value = -326.6154
$ws.Range("N55").Value = -698.45456
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 8466.134
$ws.Range("I68").Value = 2632.6667
$ws.Range("K68").Value = 2632.6667
$ws.Range("M68").Value = -1883.6667
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 8466.134
$ws.Range("I71").Value = 2632.6667
$ws.Range("K71").Value = 13163.3335
$ws.Range("M71").Value = -9419.333500000001
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 2285.8125
$ws.Range("I93").Value = 2007.1666
$ws.Range("J93").Value = 3121.75
$ws.Range("K93").Value = 2007.1666
$ws.Range("L93").Value = 3121.75
$ws.Range("M93").Value = -759.1666
$ws.Range("N93").Value = -5617.75
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 3838.7812
$ws.Range("I122").Value = 2991.8462
$ws.Range("J122").Value = 4418.263
$ws.Range("K122").Value = 8975.5386
$ws.Range("L122").Value = 13254.789
$ws.Range("M122").Value = -6525.5386
$ws.Range("N122").Value = -18154.789
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 6003
$ws.Range("I132").Value = 3693.625
$ws.Range("K132").Value = 11080.875
$ws.Range("M132").Value = -8550.875
# Row 135 (Leve Item ID 42036)
$ws.Range("H135").Value = 53999.75
$ws.Range("J135").Value = 53999.75
$ws.Range("L135").Value = 53999.75
$ws.Range("N135").Value = -64139.75
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 9434.200000000001
$ws.Range("I136").Value = 3639.8
$ws.Range("K136").Value = 10919.4
$ws.Range("M136").Value = -8369.400000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 46 (Leve Item ID 42037)
$ws.Range("H46").Value = 47797.145
$ws.Range("J46").Value = 47797.145
$ws.Range("L46").Value = 47797.145
$ws.Range("N46").Value = -48259.145
# Row 54 (Leve Item ID 3413)
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 714.44446
$ws.Range("I107").Value = 714.44446
$ws.Range("K107").Value = 2143.33338
$ws.Range("M107").Value = -223.33338
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 310715.06
$ws.Range("I132").Value = 419371.22
$ws.Range("K132").Value = 1258113.66
$ws.Range("M132").Value = -1255583.66
# Row 134 (Leve Item ID 42037)
$ws.Range("H134").Value = 47797.145
$ws.Range("J134").Value = 47797.145
$ws.Range("L134").Value = 143391.435
$ws.Range("N134").Value = -148461.435
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 5096.5835
$ws.Range("I136").Value = 2344.2856
$ws.Range("J136").Value = 8949.799999999999
$ws.Range("K136").Value = 7032.8568
$ws.Range("L136").Value = 26849.4
$ws.Range("M136").Value = -4482.8568
$ws.Range("N136").Value = -31949.4
